$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns L (12) and M (13) are formatted with a "Text" number format (numFmtId 49, "@").
# Assigning a numeric .Value directly to such a cell stores the number as literal text
# (matching real Excel behaviour for text-formatted cells). The source rows, however,
# need genuine numeric 0s written in those columns. Temporarily switching the cell's
# NumberFormat to "General" for the duration of the write forces a real numeric value
# while the original (text) NumberFormat is restored immediately afterwards, so the
# cell's style/format stays exactly as it was.
function Set-NumericValue($addr, $val) {
    $cell = $ws.Range($addr)
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "General"
    $cell.Value = $val
    $cell.NumberFormat = $fmt
}

# Rows 153-187: column M only needed a type fix (was stored as the text "0", now a
# real numeric 0 - the displayed value does not change).
for ($r = 153; $r -le 187; $r++) {
    Set-NumericValue ("M" + $r) 0
}

# Row 188: new case count for the day (C) plus the same M type fix.
$ws.Range("C188").Value = 5
Set-NumericValue "M188" 0

# Row 189: new case count for the day (C), plus type fixes on L and M.
$ws.Range("C189").Value = 14
Set-NumericValue "L189" 0
Set-NumericValue "M189" 0

# Rows 190-195: these were previously blank placeholder rows; fill in the day's
# reported figures. Columns B, H, J and K are shared formulas and recompute on
# their own once the inputs below are written.
$dailyData = @{
    190 = @{ C = 13; D = 0; E = 0; F = 0; G = 3; I = 0; L = 0; M = 0 }
    191 = @{ C = 16; D = 0; E = 0; F = 0; G = 3; I = 0; L = 0; M = 0 }
    192 = @{ C = 12; D = 0; E = 0; F = 0; G = 3; I = 0; L = 0; M = 0 }
    193 = @{ C = 9;  D = 0; E = 0; F = 0; G = 3; I = 0; L = 0; M = 0 }
    194 = @{ C = 2;  D = 0; E = 0; F = 0; G = 3; I = 0; L = 0; M = 0 }
    195 = @{ C = 0;  D = 0; E = 0; F = 0; G = 3; I = 0; L = 0; M = 0 }
}

foreach ($r in 190..195) {
    $vals = $dailyData[$r]
    $ws.Range("C" + $r).Value = $vals.C
    $ws.Range("D" + $r).Value = $vals.D
    $ws.Range("E" + $r).Value = $vals.E
    $ws.Range("F" + $r).Value = $vals.F
    $ws.Range("G" + $r).Value = $vals.G
    $ws.Range("I" + $r).Value = $vals.I
    Set-NumericValue ("L" + $r) $vals.L
    Set-NumericValue ("M" + $r) $vals.M
}

# Update the frozen-pane scroll position and the active selection to match the
# author's final view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$ws.Range("I11").Select()
